$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:G$r1")
    $rng2 = $ws.Range("B$r2`:G$r2")
    $tmp = $rng1.Value2
    $rng1.Value2 = $rng2.Value2
    $rng2.Value2 = $tmp
}

function Rotate-Rows($rows) {
    # new[rows[i]] = old[rows[i-1]]  (each row gets the content of the previous
    # row in the list; the first row gets the content of the last row)
    $vals = @()
    foreach ($r in $rows) {
        $vals += ,($ws.Range("B$r`:G$r").Value2)
    }
    $n = $rows.Count
    for ($i = 0; $i -lt $n; $i++) {
        $src = ($i - 1 + $n) % $n
        $ws.Range("B$($rows[$i])`:G$($rows[$i])").Value2 = $vals[$src]
    }
}

Swap-Rows 149 150
Swap-Rows 279 280
Swap-Rows 313 314
Swap-Rows 316 317
Swap-Rows 351 352
Swap-Rows 372 373
Swap-Rows 379 380
Swap-Rows 421 422
Swap-Rows 431 432
Swap-Rows 591 592
Swap-Rows 611 612
Swap-Rows 613 614
Swap-Rows 616 617
Swap-Rows 699 700
Swap-Rows 732 733
Swap-Rows 871 872
Swap-Rows 899 900

Rotate-Rows @(161, 162, 163)
